# Aplica a alteração: "mostrar feedback se for diferente de 00000-00-00"
# Atualiza o titulo do relatorio, os valores de recurso/glosa/banco da linha 3,
# limpa o feedback (passa a ficar em branco) e realinha a coluna de DATA.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Titulo (A1): mes de pagamento passa de JANEIRO para SETEMBRO
$ws.Range("A1").Value = "COMPETÊNCIA DE 21/10/2024 a 20/11/2024 - PGTO EM SETEMBRO 2025"

# 2) Linha 3 (dados): recurso de 200,00 concedido, reduzindo o valor a receber do banco
$ws.Range("G3").Value = 200.0     # RECURSO
$ws.Range("H3").Value = 793.85    # BANCO (993.85 - 200.00)
$ws.Range("J3").Value = 200.0     # GLOSA

# Coluna DATA (I3) passa a ficar centralizada, igual as demais colunas de data
$ws.Range("I3").HorizontalAlignment = -4108   # xlCenter

# Coluna FEEDBACK (L3) só deve exibir texto quando a data for diferente de
# 00000-00-00; nesse caso fica em branco
$ws.Range("L3").Value = ""

# 3) Ajuste fino da largura de colunas (auto-fit refletindo o novo conteudo)
$ws.Columns("G:G").ColumnWidth = 13.166666666666666
$ws.Columns("I:I").ColumnWidth = 14.333333333333334
$ws.Columns("J:J").ColumnWidth = 15.5
$ws.Columns("L:L").ColumnWidth = 13.166666666666666
